$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Fixed error in input calculations"
# The formula in B6 was computing B4/B5 (costPV / capPV) which is wrong;
# it should compute B5/B3 (capPV / etaPV_rated).
$ws.Range("B6").Formula = "=B5/B3"

# Move the active selection on the sheet from B9 to G9, matching the
# cursor position recorded in the saved workbook.
$ws.Range("G9").Select()
